# Update the crypto listing sheet with refreshed prices/volumes (and a
# couple of rows whose coin identity shifted rank position), as produced
# by the scheduled GitHub Actions data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '49.987.15'
    'E2' = '  +3.85%  '
    'D3' = '2.656.02'
    'E3' = '  +6.09%  '
    'E4' = '  +0.03%  '
    'D5' = '113.84'
    'E5' = '  +7.41%  '
    'D6' = '327.91'
    'E6' = '  +2.85%  '
    'E7' = '  +1.63%  '
    'E8' = '  +0.08%  '
    'E9' = '  +3.55%  '
    'D10' = '41.31'
    'E10' = '  +6.19%  '
    'D11' = '20.18'
    'E11' = '  +1.60%  '
    'E12' = '  +2.83%  '
    'E13' = '  +1.02%  '
    'D14' = '7.36'
    'E14' = '  +4.41%  '
    'D15' = '3.008.31'
    'E15' = '  +3.83%  '
    'D16' = '2.631.13'
    'E16' = '  +4.92%  '
    'E17' = '  +5.40%  '
    'D18' = '49.936.83'
    'E18' = '  +3.94%  '
    'D19' = '13.18'
    'E19' = '  +1.95%  '
    'D20' = '6.79'
    'E20' = '  +2.36%  '
    'E21' = '  -0.15%  '
    'E22' = '  +3.08%  '
    'D23' = '72.42'
    'E23' = '  +2.04%  '
    'D24' = '277.79'
    'E24' = '  +2.10%  '
    'E25' = '  +3.46%  '
    'D26' = '26.88'
    'E26' = '  +4.18%  '
    'E27' = '  -0.03%  '
    'D28' = '9.96'
    'E28' = '  +2.73%  '
    'E29' = '  +1.13%  '
    'D30' = '36.33'
    'E30' = '  +5.14%  '
    'E31' = '  -2.06%  '
    'D32' = '50.34'
    'E32' = '  +2.05%  '
    'B33' = 'Celestia'
    'C33' = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
    'D33' = '19.57'
    'E33' = '  +2.49%  '
    'B34' = 'Filecoin'
    'C34' = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
    'D34' = '5.43'
    'E34' = '  +2.78%  '
    'B35' = 'Hedera'
    'C35' = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
    'D35' = '0.0802'
    'E35' = '  +3.52%  '
    'B36' = 'FirstDigitalUSD'
    'C36' = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
    'D36' = '1.00'
    'E36' = '  -0.24%  '
    'E37' = '  +7.15%  '
    'E38' = '  +4.16%  '
    'E39' = '  +7.60%  '
    'D40' = '126.01'
    'E40' = '  +5.07%  '
    'E41' = '  +1.91%  '
    'B42' = 'WEMIXToken'
    'C42' = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
    'D42' = '2.24'
    'E42' = '  +1.20%  '
    'B43' = 'EnergySwap'
    'C43' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'D43' = '22.43'
    'E43' = '  +3.32%  '
    'E44' = '  +4.13%  '
    'E45' = '  +5.38%  '
    'D46' = '2.073.12'
    'E46' = '  +3.62%  '
    'E47' = '  +12.46%  '
    'D48' = '1.99'
    'E48' = '  +6.13%  '
    'E49' = '  +2.61%  '
    'D50' = '5.42'
    'E50' = '  +4.62%  '
    'D51' = '82.07'
    'E51' = '  +4.22%  '
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    # Force text interpretation (so values like "1.00" or "113.84" are not
    # auto-converted to numbers), then strip the number-format override we
    # just applied so the cell's style stays exactly as it was originally
    # (unstyled / default), matching the source data which stores these as
    # plain text cells.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.ClearFormats()
}
